# Update "想去人数" (number of people interested) counts on three sheets:
#   Sheet 1 "展览"     (Exhibitions)
#   Sheet 2 "演出"     (Performances)
#   Sheet 4 "全部类型" (All types)
# Sheet 3 "本地生活" has no changes.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value  = 46
$ws1.Range("F8").Value  = 6882
$ws1.Range("F12").Value = 6456
$ws1.Range("F15").Value = 4379
$ws1.Range("F19").Value = 4374
$ws1.Range("F21").Value = 236
$ws1.Range("F22").Value = 238
$ws1.Range("F23").Value = 325
$ws1.Range("F29").Value = 73
$ws1.Range("F30").Value = 7928
$ws1.Range("F38").Value = 1601
$ws1.Range("F39").Value = 206
$ws1.Range("F40").Value = 918
$ws1.Range("F42").Value = 3999
$ws1.Range("F46").Value = 41
$ws1.Range("F49").Value = 6

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value  = 21
$ws2.Range("F11").Value = 159
$ws2.Range("F19").Value = 871

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value  = 21
$ws4.Range("F8").Value  = 46
$ws4.Range("F12").Value = 6882
$ws4.Range("F16").Value = 6456
$ws4.Range("F19").Value = 4379
$ws4.Range("F22").Value = 4374
$ws4.Range("F24").Value = 236
$ws4.Range("F25").Value = 238
$ws4.Range("F26").Value = 325
$ws4.Range("F29").Value = 159
$ws4.Range("F30").Value = 7928
$ws4.Range("F38").Value = 1601
$ws4.Range("F39").Value = 206
$ws4.Range("F40").Value = 918
$ws4.Range("F42").Value = 4000
$ws4.Range("F49").Value = 6
